$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 783.125
$ws.Range("I111").Value = 899.25
$ws.Range("K111").Value = 2697.75
$ws.Range("M111").Value = 369.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1156.6897
$ws.Range("I2").Value = 1024.6666
$ws.Range("K2").Value = 1024.6666
$ws.Range("M2").Value = -911.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2888.5
$ws.Range("I32").Value = 2908.725
$ws.Range("K32").Value = 2908.725
$ws.Range("M32").Value = -2621.725

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1156.6897
$ws.Range("I116").Value = 1024.6666
$ws.Range("K116").Value = 1024.6666
$ws.Range("M116").Value = 1269.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2925
$ws.Range("I132").Value = 1353.3334
$ws.Range("K132").Value = 4060.0002
$ws.Range("M132").Value = -1530.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1156.6897
$ws.Range("I3").Value = 1024.6666
$ws.Range("K3").Value = 1024.6666
$ws.Range("M3").Value = -910.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 997.74194
$ws.Range("I134").Value = 790.64
$ws.Range("K134").Value = 2371.92
$ws.Range("M134").Value = 163.0799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 75000

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 40878.145
$ws.Range("I31").Value = 42560.16
$ws.Range("J31").Value = 36673.1
$ws.Range("K31").Value = 42560.16
$ws.Range("L31").Value = 36673.1
$ws.Range("M31").Value = -42265.16
$ws.Range("N31").Value = -37263.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 40878.145
$ws.Range("I34").Value = 42560.16
$ws.Range("J34").Value = 36673.1
$ws.Range("K34").Value = 42560.16
$ws.Range("L34").Value = 36673.1
$ws.Range("M34").Value = -42358.16
$ws.Range("N34").Value = -37077.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 15580.667
$ws.Range("I62").Value = 17297.4
$ws.Range("K62").Value = 17297.4
$ws.Range("M62").Value = -16673.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 15580.667
$ws.Range("I65").Value = 17297.4
$ws.Range("K65").Value = 86487
$ws.Range("M65").Value = -83367

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4044.5
$ws.Range("I132").Value = 3714.0667
$ws.Range("J132").Value = 5696.6665
$ws.Range("K132").Value = 11142.2001
$ws.Range("L132").Value = 17089.9995
$ws.Range("M132").Value = -8612.2001
$ws.Range("N132").Value = -22149.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 755.54285
$ws.Range("I5").Value = 705.65515
$ws.Range("K5").Value = 2116.96545
$ws.Range("M5").Value = -2004.96545

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 93.95
$ws.Range("J12").Value = 103.86667
$ws.Range("L12").Value = 311.60001
$ws.Range("N12").Value = -657.60001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 8363.904
$ws.Range("I56").Value = 8363.904
$ws.Range("K56").Value = 8363.904
$ws.Range("M56").Value = -7833.904

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 2989.1
$ws.Range("I64").Value = 1978.4
$ws.Range("K64").Value = 5935.200000000001
$ws.Range("M64").Value = -5665.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 2989.1
$ws.Range("I67").Value = 1978.4
$ws.Range("K67").Value = 5935.200000000001
$ws.Range("M67").Value = -4999.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 760
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 700
$ws.Range("K107").Value = 3000
$ws.Range("L107").Value = 2100
$ws.Range("M107").Value = -1080
$ws.Range("N107").Value = -5940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2129.8
$ws.Range("J113").Value = 2162.25
$ws.Range("L113").Value = 6486.75
$ws.Range("N113").Value = -10826.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 755.54285
$ws.Range("I135").Value = 705.65515
$ws.Range("K135").Value = 6350.896350000001
$ws.Range("M135").Value = -3815.896350000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 29999
$ws.Range("J57").Value = 29999
$ws.Range("L57").Value = 29999
$ws.Range("N57").Value = -31639

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3034.7646
$ws.Range("I80").Value = 2711.261
$ws.Range("K80").Value = 2711.261
$ws.Range("M80").Value = -1713.261

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3034.7646
$ws.Range("I83").Value = 2711.261
$ws.Range("K83").Value = 13556.305
$ws.Range("M83").Value = -8564.305

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4415.3887
$ws.Range("I132").Value = 3531.561
$ws.Range("J132").Value = 7202.846
$ws.Range("K132").Value = 10594.683
$ws.Range("L132").Value = 21608.538
$ws.Range("M132").Value = -8064.683000000001
$ws.Range("N132").Value = -26668.538

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 118007
$ws.Range("J140").Value = 118007
$ws.Range("L140").Value = 118007
$ws.Range("N140").Value = -128367

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11154.171
$ws.Range("I7").Value = 18737.895
$ws.Range("K7").Value = 18737.895
$ws.Range("M7").Value = -18625.895

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 999.5
$ws.Range("J16").Value = 999
$ws.Range("L16").Value = 999
$ws.Range("N16").Value = -1339

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5609.8887
$ws.Range("I40").Value = 4372.875
$ws.Range("K40").Value = 4372.875
$ws.Range("M40").Value = -4236.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 32828.07
$ws.Range("I46").Value = 61699
$ws.Range("J46").Value = 3957.1428
$ws.Range("K46").Value = 61699
$ws.Range("L46").Value = 3957.1428
$ws.Range("M46").Value = -61511
$ws.Range("N46").Value = -4333.1428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3622
$ws.Range("I82").Value = 3622
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 3622
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -3261
$ws.Range("N82").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 3622
$ws.Range("I85").Value = 3622
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 3622
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -2374
$ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 11154.171
$ws.Range("I126").Value = 18737.895
$ws.Range("K126").Value = 56213.685
$ws.Range("M126").Value = -53743.685

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 67399.39999999999
$ws.Range("J133").Value = 67399.39999999999
$ws.Range("L133").Value = 67399.39999999999
$ws.Range("N133").Value = -77519.39999999999
